$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 19, pushing existing rows (19..34) down by one.
# This mirrors row 19's existing content into the freshly inserted row, then
# the date (D) and volume (J) get their new values below.
$ws.Rows.Item(19).Insert()

$newRow = 19
$srcRow = 20

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($col in $cols) {
    $src = $ws.Range("$col$srcRow")
    $dst = $ws.Range("$col$newRow")
    $dst.Value = $src.Value()
}
$ws.Range("D$newRow").NumberFormat = $ws.Range("D$srcRow").NumberFormat()

# Apply the new values for the inserted row (date & volume).
$ws.Range("D19").Value = 44874
$ws.Range("J19").Value = 7900

# Update values that changed further down the sheet because of the shift.
$ws.Range("K32").Value = 3000
$ws.Range("M32").Value = 3000
$ws.Range("P32").Value = 30

$ws.Range("J33").Value = 7000
$ws.Range("K33").Value = 2500
$ws.Range("M33").Value = 2750
$ws.Range("P33").Value = 28
